# Applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions" — updates the Price (D)
# and Volume(1h) (E) columns for each coin row, and fixes two pairs of
# rows (12/13 and 19/20) whose Coin name + Link got swapped back to the
# correct coin.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.583.83'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.921.32'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.78%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4728'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2900'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06787'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '105.04'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '18.34'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.910.35'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.26%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07709'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.299'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6738'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '287.15'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.26%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.608.16'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007628'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.33%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.94'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.163.64'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.435'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.321'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.385'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.44'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.80'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.121'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +10.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1078'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('E30').Value = '  +3.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.192'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.130'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05042'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7431'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.154'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02073'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.744'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.694'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.061'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '111.17'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8810'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4371'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.892'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '67.13'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.244'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.252'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '47.90'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +16.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1233'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.91'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4045'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.42%  '
